# Generate Report for Handoff
#
# Flips the localization-status workbook from "handed back" to "ready for
# handoff": updates the Status text wherever it appears (Overview zh-cn/de-de
# columns, and the Status column on each language sheet), bumps the
# "generate" timestamps that go with the new handoff, and narrows the
# Status-ish columns that used to be sized for the long "Handed back..."
# string.

$wb = $excel.ActiveWorkbook

$newStatus = "Ready for handoff"

# "Latest HO Xliff Generate Date" (Overview) / "Latest Handoff Datetime" (de-de)
$newOverviewDate = "2016-08-31 07:08:18"
# "Latest Handoff Datetime" (zh-cn)
$newZhCnDate = "2016-08-31 07:08:14"

# New (narrower) column width - matches the width used elsewhere for
# compact date/status columns.
$newColumnWidth = 16.333333333333332

# --- Overview sheet ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = $newStatus
$wsOverview.Range("F2").Value = $newStatus
$wsOverview.Range("G2").Value = $newOverviewDate
$wsOverview.Range("E1:F1").ColumnWidth = $newColumnWidth

# --- zh-cn sheet ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = $newStatus
$wsZhCn.Range("H2").Value = $newZhCnDate
$wsZhCn.Range("C1").ColumnWidth = $newColumnWidth

# --- de-de sheet ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = $newStatus
$wsDeDe.Range("H2").Value = $newOverviewDate
$wsDeDe.Range("C1").ColumnWidth = $newColumnWidth
